$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Update existing rows 76-97 with new values
$values76to97 = @(79,82,81,84,83,86,85,88,87,90,89,92,48,33,50,27,36,45,39,28,34,41)
$row = 76
foreach ($v in $values76to97) {
    $ws.Cells.Item($row, 1).Value = $v
    $row++
}

# Add new rows 98-100
$ws.Cells.Item(98, 1).Value = 30
$ws.Cells.Item(99, 1).Value = 38
$ws.Cells.Item(100, 1).Value = 44

# Update selection to match target state
$ws.Range("A101:A103").Select()
